# Daily attendance processing - 2025-11-22 07:20:20
#
# The "Recorded By" column (column G) lists the accounts that touched each
# attendance record as a comma separated string, e.g.
#   "dnasr281@gmail.com, System"
# For this processing pass the last two recorders in that list need to be
# swapped (the two most recent entries trade places), e.g.
#   "dnasr281@gmail.com, System"  ->  "System, dnasr281@gmail.com"
#   "System, admin@admin.com"     ->  "admin@admin.com, System"
# Cells that only contain a single recorder are left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$used = $ws.UsedRange
$firstRow = $used.Row
$lastRow = $firstRow + $used.Rows.Count - 1

# Column G = "Recorded By"
$col = 7

for ($r = $firstRow; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, $col)
    $val = $cell.Value()

    if ($val -eq $null) {
        continue
    }

    $sval = [string]$val

    if ($sval.Contains(",")) {
        $parts = $sval.Split(",")
        $n = $parts.Length

        for ($i = 0; $i -lt $n; $i++) {
            $parts[$i] = $parts[$i].Trim()
        }

        if ($n -ge 2) {
            $tmp = $parts[$n - 1]
            $parts[$n - 1] = $parts[$n - 2]
            $parts[$n - 2] = $tmp

            $newval = [string]::Join(", ", $parts)

            if ($newval -ne $sval) {
                $cell.Value = $newval
            }
        }
    }
}
